# Update TPM-derived values in the LR-pairs sheet (Inhba-Acvr1)
# Only numeric result columns (G..T) change; identifiers/columns A-F,K,L stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ G=0.08097566666666667; H=0.242927; I=0.005588990034505014; J=0.005588990034505015;
            M=4.621579; N=13.864737; O=0.1778708528171788; P=0.1778708528171788;
            Q=0.3742354405776667; R=3.368118965199; S=0.0009941184238241203; T=0.0009941184238241205 }
    3  = @{ G=0.08097566666666667; H=0.242927; I=0.005588990034505014; J=0.005588990034505015;
            N=46.543441; O=0.5971062807549863; P=0.5971062807549863;
            Q=1.256295387978556; R=11.306658491807; S=0.003337221052679971; T=0.003337221052679972 }
    4  = @{ G=0.08097566666666667; H=0.242927; I=0.005588990034505014; J=0.005588990034505015;
            O=0.2250228664278349; P=0.2250228664278349;
            Q=0.4734419958295555; R=4.260977962466; S=0.001257650558000922; T=0.001257650558000923 }
    5  = @{ I=0.6976944377922635; J=0.6976944377922635;
            M=4.621579; N=13.864737; O=0.1778708528171788; P=0.1778708528171788;
            Q=46.717203592026; R=420.454832328234; S=0.124099504655912; T=0.124099504655912 }
    6  = @{ I=0.6976944377922635; J=0.6976944377922635;
            N=46.543441; O=0.5971062807549863; P=0.5971062807549863;
            S=0.4165977308535796; T=0.4165977308535796 }
    7  = @{ I=0.6976944377922635; J=0.6976944377922635;
            O=0.2250228664278349; P=0.2250228664278349;
            S=0.1569972022827719; T=0.1569972022827719 }
    8  = @{ I=0.2967165721732315; J=0.2967165721732316;
            M=4.621579; N=13.864737; O=0.1778708528171788; P=0.1778708528171788;
            Q=19.867964771524; R=178.811682943716; S=0.05277722973744267; T=0.05277722973744269 }
    9  = @{ I=0.2967165721732315; J=0.2967165721732316;
            N=46.543441; O=0.5971062807549863; P=0.5971062807549863;
            Q=66.69606831586535; R=600.2646148427881; S=0.1771713288487267; T=0.1771713288487267 }
    10 = @{ I=0.2967165721732315; J=0.2967165721732316;
            O=0.2250228664278349; P=0.2250228664278349;
            S=0.06676801358706212; T=0.06676801358706214 }
}

foreach ($rowKey in $data.Keys) {
    $rowNum = [int]$rowKey
    $cols = $data[$rowKey]
    foreach ($colKey in $cols.Keys) {
        $cellAddr = "$colKey$rowNum"
        $ws.Range($cellAddr).Value = $cols[$colKey]
    }
}
